$d = $word.ActiveDocument

$d.Content.Find.Execute("91-49=42", $true, $false, $false, $false, $false, $true, 1, $false, "35-19=16", 2) | Out-Null
$d.Content.Find.Execute("65-16=49", $true, $false, $false, $false, $false, $true, 1, $false, "25-23=2", 2) | Out-Null
$d.Content.Find.Execute("83-77=6", $true, $false, $false, $false, $false, $true, 1, $false, "50-34=16", 2) | Out-Null
$d.Content.Find.Execute("74-21=53", $true, $false, $false, $false, $false, $true, 1, $false, "72+0=72", 2) | Out-Null
$d.Content.Find.Execute("21+52=73", $true, $false, $false, $false, $false, $true, 1, $false, "8+43=51", 2) | Out-Null
$d.Content.Find.Execute("8+81=89", $true, $false, $false, $false, $false, $true, 1, $false, "35+53=88", 2) | Out-Null
$d.Content.Find.Execute("20+51=71", $true, $false, $false, $false, $false, $true, 1, $false, "59+4=63", 2) | Out-Null
$d.Content.Find.Execute("85-66=19", $true, $false, $false, $false, $false, $true, 1, $false, "54-43=11", 2) | Out-Null
$d.Content.Find.Execute("34-14=20", $true, $false, $false, $false, $false, $true, 1, $false, "89-64=25", 2) | Out-Null
$d.Content.Find.Execute("53+6=59", $true, $false, $false, $false, $false, $true, 1, $false, "96-3=93", 2) | Out-Null
$d.Content.Find.Execute("45-20=25", $true, $false, $false, $false, $false, $true, 1, $false, "75-66=9", 2) | Out-Null
$d.Content.Find.Execute("11+32=43", $true, $false, $false, $false, $false, $true, 1, $false, "15+57=72", 2) | Out-Null
$d.Content.Find.Execute("3+89=92", $true, $false, $false, $false, $false, $true, 1, $false, "24+51=75", 2) | Out-Null
$d.Content.Find.Execute("85-65=20", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=22", 2) | Out-Null
$d.Content.Find.Execute("71+17=88", $true, $false, $false, $false, $false, $true, 1, $false, "25+45=70", 2) | Out-Null
$d.Content.Find.Execute("84+8=92", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=4", 2) | Out-Null
$d.Content.Find.Execute("87-27=60", $true, $false, $false, $false, $false, $true, 1, $false, "31-3=28", 2) | Out-Null
$d.Content.Find.Execute("44-22=22", $true, $false, $false, $false, $false, $true, 1, $false, "13+2=15", 2) | Out-Null
$d.Content.Find.Execute("90-20=70", $true, $false, $false, $false, $false, $true, 1, $false, "54-6=48", 2) | Out-Null
$d.Content.Find.Execute("96-27=69", $true, $false, $false, $false, $false, $true, 1, $false, "42-4=38", 2) | Out-Null
$d.Content.Find.Execute("54-48=6", $true, $false, $false, $false, $false, $true, 1, $false, "79-65=14", 2) | Out-Null
$d.Content.Find.Execute("2+26=28", $true, $false, $false, $false, $false, $true, 1, $false, "24-21=3", 2) | Out-Null
$d.Content.Find.Execute("92-34=58", $true, $false, $false, $false, $false, $true, 1, $false, "66-52=14", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=65", 2) | Out-Null
$d.Content.Find.Execute("19-4=15", $true, $false, $false, $false, $false, $true, 1, $false, "4+28=32", 2) | Out-Null
$d.Content.Find.Execute("48+3=51", $true, $false, $false, $false, $false, $true, 1, $false, "88-26=62", 2) | Out-Null
$d.Content.Find.Execute("45+18=63", $true, $false, $false, $false, $false, $true, 1, $false, "93-74=19", 2) | Out-Null
$d.Content.Find.Execute("10+47=57", $true, $false, $false, $false, $false, $true, 1, $false, "83-64=19", 2) | Out-Null
$d.Content.Find.Execute("21+11=32", $true, $false, $false, $false, $false, $true, 1, $false, "26+71=97", 2) | Out-Null
$d.Content.Find.Execute("29-1=28", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=43", 2) | Out-Null
$d.Content.Find.Execute("14+63=77", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=59", 2) | Out-Null
$d.Content.Find.Execute("76-57=19", $true, $false, $false, $false, $false, $true, 1, $false, "86-50=36", 2) | Out-Null
$d.Content.Find.Execute("59-51=8", $true, $false, $false, $false, $false, $true, 1, $false, "52-38=14", 2) | Out-Null
$d.Content.Find.Execute("27+20=47", $true, $false, $false, $false, $false, $true, 1, $false, "27+58=85", 2) | Out-Null
$d.Content.Find.Execute("24+9=33", $true, $false, $false, $false, $false, $true, 1, $false, "20+36=56", 2) | Out-Null
$d.Content.Find.Execute("13+36=49", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=9", 2) | Out-Null
$d.Content.Find.Execute("44-0=44", $true, $false, $false, $false, $false, $true, 1, $false, "69-29=40", 2) | Out-Null
$d.Content.Find.Execute("17+33=50", $true, $false, $false, $false, $false, $true, 1, $false, "69-26=43", 2) | Out-Null
$d.Content.Find.Execute("47-40=7", $true, $false, $false, $false, $false, $true, 1, $false, "25+32=57", 2) | Out-Null
$d.Content.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=16", 2) | Out-Null
$d.Content.Find.Execute("86-29=57", $true, $false, $false, $false, $false, $true, 1, $false, "59-57=2", 2) | Out-Null
$d.Content.Find.Execute("48+1=49", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=75", 2) | Out-Null
$d.Content.Find.Execute("2+97=99", $true, $false, $false, $false, $false, $true, 1, $false, "6+38=44", 2) | Out-Null
$d.Content.Find.Execute("9+11=20", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=48", 2) | Out-Null
$d.Content.Find.Execute("35+51=86", $true, $false, $false, $false, $false, $true, 1, $false, "33+14=47", 2) | Out-Null
$d.Content.Find.Execute("75+14=89", $true, $false, $false, $false, $false, $true, 1, $false, "12+9=21", 2) | Out-Null
$d.Content.Find.Execute("13+66=79", $true, $false, $false, $false, $false, $true, 1, $false, "25-24=1", 2) | Out-Null
$d.Content.Find.Execute("57+24=81", $true, $false, $false, $false, $false, $true, 1, $false, "74-16=58", 2) | Out-Null
$d.Content.Find.Execute("60+6=66", $true, $false, $false, $false, $false, $true, 1, $false, "61-21=40", 2) | Out-Null
$d.Content.Find.Execute("67-11=56", $true, $false, $false, $false, $false, $true, 1, $false, "52-26=26", 2) | Out-Null
$d.Content.Find.Execute("78-20=58", $true, $false, $false, $false, $false, $true, 1, $false, "25-1=24", 2) | Out-Null
$d.Content.Find.Execute("45+22=67", $true, $false, $false, $false, $false, $true, 1, $false, "17+29=46", 2) | Out-Null
$d.Content.Find.Execute("69+6=75", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=81", 2) | Out-Null
$d.Content.Find.Execute("34-5=29", $true, $false, $false, $false, $false, $true, 1, $false, "56-30=26", 2) | Out-Null
$d.Content.Find.Execute("76-13=63", $true, $false, $false, $false, $false, $true, 1, $false, "46+28=74", 2) | Out-Null
$d.Content.Find.Execute("78-32=46", $true, $false, $false, $false, $false, $true, 1, $false, "92-82=10", 2) | Out-Null
$d.Content.Find.Execute("47-46=1", $true, $false, $false, $false, $false, $true, 1, $false, "52+41=93", 2) | Out-Null
$d.Content.Find.Execute("78-39=39", $true, $false, $false, $false, $false, $true, 1, $false, "73-27=46", 2) | Out-Null
$d.Content.Find.Execute("23+8=31", $true, $false, $false, $false, $false, $true, 1, $false, "9+51=60", 2) | Out-Null
$d.Content.Find.Execute("89-69=20", $true, $false, $false, $false, $false, $true, 1, $false, "96-1=95", 2) | Out-Null
$d.Content.Find.Execute("86-5=81", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=91", 2) | Out-Null
$d.Content.Find.Execute("21-3=18", $true, $false, $false, $false, $false, $true, 1, $false, "89-69=20", 2) | Out-Null
$d.Content.Find.Execute("56-36=20", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=46", 2) | Out-Null
$d.Content.Find.Execute("35-31=4", $true, $false, $false, $false, $false, $true, 1, $false, "50+44=94", 2) | Out-Null
$d.Content.Find.Execute("54-17=37", $true, $false, $false, $false, $false, $true, 1, $false, "96+3=99", 2) | Out-Null
$d.Content.Find.Execute("3+33=36", $true, $false, $false, $false, $false, $true, 1, $false, "36+31=67", 2) | Out-Null
$d.Content.Find.Execute("62-10=52", $true, $false, $false, $false, $false, $true, 1, $false, "92-24=68", 2) | Out-Null
$d.Content.Find.Execute("84-50=34", $true, $false, $false, $false, $false, $true, 1, $false, "89+3=92", 2) | Out-Null
$d.Content.Find.Execute("78+13=91", $true, $false, $false, $false, $false, $true, 1, $false, "51+32=83", 2) | Out-Null
$d.Content.Find.Execute("2+61=63", $true, $false, $false, $false, $false, $true, 1, $false, "40-20=20", 2) | Out-Null
$d.Content.Find.Execute("32+8=40", $true, $false, $false, $false, $false, $true, 1, $false, "87-22=65", 2) | Out-Null
$d.Content.Find.Execute("26+23=49", $true, $false, $false, $false, $false, $true, 1, $false, "52+25=77", 2) | Out-Null
$d.Content.Find.Execute("94-1=93", $true, $false, $false, $false, $false, $true, 1, $false, "94-89=5", 2) | Out-Null
$d.Content.Find.Execute("78-34=44", $true, $false, $false, $false, $false, $true, 1, $false, "23+31=54", 2) | Out-Null
$d.Content.Find.Execute("92-66=26", $true, $false, $false, $false, $false, $true, 1, $false, "69+25=94", 2) | Out-Null
$d.Content.Find.Execute("35+6=41", $true, $false, $false, $false, $false, $true, 1, $false, "68+5=73", 2) | Out-Null
$d.Content.Find.Execute("2+22=24", $true, $false, $false, $false, $false, $true, 1, $false, "76-28=48", 2) | Out-Null
$d.Content.Find.Execute("36-18=18", $true, $false, $false, $false, $false, $true, 1, $false, "78-44=34", 2) | Out-Null
$d.Content.Find.Execute("47+18=65", $true, $false, $false, $false, $false, $true, 1, $false, "6+19=25", 2) | Out-Null
$d.Content.Find.Execute("60-45=15", $true, $false, $false, $false, $false, $true, 1, $false, "82+3=85", 2) | Out-Null
$d.Content.Find.Execute("90-11=79", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=57", 2) | Out-Null
$d.Content.Find.Execute("36+24=60", $true, $false, $false, $false, $false, $true, 1, $false, "53-20=33", 2) | Out-Null
$d.Content.Find.Execute("0+43=43", $true, $false, $false, $false, $false, $true, 1, $false, "13+27=40", 2) | Out-Null
$d.Content.Find.Execute("71-8=63", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=55", 2) | Out-Null
$d.Content.Find.Execute("25+56=81", $true, $false, $false, $false, $false, $true, 1, $false, "76-34=42", 2) | Out-Null
$d.Content.Find.Execute("95-87=8", $true, $false, $false, $false, $false, $true, 1, $false, "85+12=97", 2) | Out-Null
$d.Content.Find.Execute("76+21=97", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=53", 2) | Out-Null
$d.Content.Find.Execute("17+11=28", $true, $false, $false, $false, $false, $true, 1, $false, "87+6=93", 2) | Out-Null
$d.Content.Find.Execute("1+83=84", $true, $false, $false, $false, $false, $true, 1, $false, "68+8=76", 2) | Out-Null
$d.Content.Find.Execute("84-35=49", $true, $false, $false, $false, $false, $true, 1, $false, "9+29=38", 2) | Out-Null
$d.Content.Find.Execute("58+5=63", $true, $false, $false, $false, $false, $true, 1, $false, "8+32=40", 2) | Out-Null
$d.Content.Find.Execute("75-27=48", $true, $false, $false, $false, $false, $true, 1, $false, "14+38=52", 2) | Out-Null
$d.Content.Find.Execute("24+56=80", $true, $false, $false, $false, $false, $true, 1, $false, "66-3=63", 2) | Out-Null
$d.Content.Find.Execute("34-20=14", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=45", 2) | Out-Null
$d.Content.Find.Execute("61-18=43", $true, $false, $false, $false, $false, $true, 1, $false, "77-72=5", 2) | Out-Null
$d.Content.Find.Execute("26+58=84", $true, $false, $false, $false, $false, $true, 1, $false, "83-35=48", 2) | Out-Null
$d.Content.Find.Execute("61-48=13", $true, $false, $false, $false, $false, $true, 1, $false, "45-26=19", 2) | Out-Null
$d.Content.Find.Execute("37+20=57", $true, $false, $false, $false, $false, $true, 1, $false, "72-0=72", 2) | Out-Null
$d.Content.Find.Execute("70+29=99", $true, $false, $false, $false, $false, $true, 1, $false, "89-64=25", 2) | Out-Null
$d.Content.Find.Execute("72-40=32", $true, $false, $false, $false, $false, $true, 1, $false, "99-0=99", 2) | Out-Null
